$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "总计" (summary) sheet: insert a new data row for 2022-Q3 right after the
#    header row, pushing the existing quarters down by one row.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$summary.Rows.Item(2).Insert()

# Copy the (now shifted-down) formatting of the old row 2 onto the freshly
# inserted blank row so the new row's styles (bold index col, etc.) line up
# with the rest of the table.
$summary.Range("A3:D3").Copy()
$summary.Range("A2:D2").PasteSpecial(-4122)

$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "2022-Q3"
$summary.Cells.Item(2, 3).Value = 1
$summary.Cells.Item(2, 4).Value = 0.08

# ---------------------------------------------------------------------------
# 2) New "2022-Q3" sheet: duplicate the "2022-Q2" sheet (same column layout)
#    and place it right before it, then overwrite the data row with the new
#    quarter's numbers.
# ---------------------------------------------------------------------------
$prevQ = $wb.Worksheets.Item("2022-Q2")
$prevQ.Copy($prevQ)
$newQ = $wb.Worksheets.Item("2022-Q2 (2)")
$newQ.Name = "2022-Q3"

$newQ.Range("C2").Value = "工银全球精选股票（QDII）"
$newQ.Range("D2").Value = "'3.72"
$newQ.Range("E2").Value = "'93.69"
$newQ.Range("F2").Value = "'2.24"
$newQ.Range("G2").Value = "'0.0833"
$newQ.Range("H2").Value = 5

# Restore the originally-selected tab (last quarter sheet) as the active one,
# since copying a sheet makes the copy the active tab.
$wb.Worksheets.Item("2020-Q4").Activate()
